$d = $word.ActiveDocument

# --- Merge split runs back into single runs (no textual change, just
#     collapsing the run boundaries that existed in the "before" doc) ---

$null = $d.Content.Find.Execute(
    "We aim to learn how to collaborate as a team to effectively develop software.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "We aim to learn how to collaborate as a team to effectively develop software.", 2)

$null = $d.Content.Find.Execute(
    " attendance at meetings, participation, frequency of communication, the quality of work, etc.?",
    $false, $false, $false, $false, $false, $true, 1, $false,
    " attendance at meetings, participation, frequency of communication, the quality of work, etc.?", 2)

$null = $d.Content.Find.Execute(
    "we expect to draft the first version in the first 3 weeks to refine it before the deadline.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "we expect to draft the first version in the first 3 weeks to refine it before the deadline.", 2)

$null = $d.Content.Find.Execute(
    "ROLES: Which roles do we need in this project and how do we allocate them? Will there be a project lead?  Is there a need to rotate roles?",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "ROLES: Which roles do we need in this project and how do we allocate them? Will there be a project lead?  Is there a need to rotate roles?", 2)

$rsquo = [char]8217
$null = $d.Content.Find.Execute(
    "We don" + $rsquo + "t need fixed roles.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "We don" + $rsquo + "t need fixed roles.", 2)

$null = $d.Content.Find.Execute(
    "CONSEQUENCES: How will we address non-performance regarding these goals, expectations, policies and procedures?  How do we resolve disagreements?",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "CONSEQUENCES: How will we address non-performance regarding these goals, expectations, policies and procedures?  How do we resolve disagreements?", 2)

$null = $d.Content.Find.Execute(
    "If a team member does not deliver on time, another team member will cover.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "If a team member does not deliver on time, another team member will cover.", 2)

# --- Sign the contract: fill in the two still-blank signature lines ---

$pavendran = $d.Paragraphs(51)
$rngP = $pavendran.Range
$rngP.InsertBefore("Pavendran Wimalendran 14/March/2022")
$rngPFmt = $pavendran.Range.Font
$rngPFmt.NameAscii = "Calibri"
$rngPFmt.NameFarEast = "Calibri"
$rngPFmt.NameOther = "Calibri"
$rngPFmt.NameBi = "Calibri"

$victoria = $d.Paragraphs(54)
$rngV = $victoria.Range
$rngV.InsertBefore("Victoria Thompson 15/March/2022")
$rngVFmt = $victoria.Range.Font
$rngVFmt.NameAscii = "Calibri"
$rngVFmt.NameFarEast = "Calibri"
$rngVFmt.NameOther = "Calibri"
$rngVFmt.NameBi = "Calibri"
